# Generate Report for Handoff
# Updates the "Status" of the localization job from "In Translation" to
# "Ready for handoff" across the Overview / zh-cn / de-de sheets, and
# refreshes the associated "Latest HO Xliff Generate Date" / "Latest
# Handoff Datetime" timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-08-27 02:57:08"

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-27 02:56:59"

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-27 02:57:08"

# The longer "Ready for handoff" status text no longer fits the old
# column width, so the status columns widen (mirrors Excel's own
# autofit-on-edit behaviour for these generated reports).
$ws = $wb.Worksheets.Item("Overview")
$ws.Columns.Item(5).ColumnWidth = 16.333333333333332
$ws.Columns.Item(6).ColumnWidth = 16.333333333333332

$ws = $wb.Worksheets.Item("zh-cn")
$ws.Columns.Item(3).ColumnWidth = 16.333333333333332

$ws = $wb.Worksheets.Item("de-de")
$ws.Columns.Item(3).ColumnWidth = 16.333333333333332
